# Updates the cryptos price/volume table with freshly scraped values.
# Note: several "Price" column values look like plain numbers (e.g. "1.00",
# "57.73") but must stay as literal text (matching the original inlineStr
# cells, which preserve trailing zeros / exact formatting). For those we
# prefix the value with a leading apostrophe to force Excel to store it as
# text instead of auto-converting to a float, then reset the cell style to
# "Normal" so no stray number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.740.96"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").Value = "3.669.48"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'203.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.06%  "
$ws.Range("D6").Value = "'581.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "3.665.86"
$ws.Range("E7").Value = "  +2.82%  "
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("E11").Value = "  +8.88%  "
$ws.Range("D12").Value = "'57.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.33%  "
$ws.Range("D13").Value = "'0.0000298"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +18.76%  "
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("D15").Value = "4.261.62"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "3.671.95"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "'12.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.03%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "68.790.79"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'18.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("E21").Value = "  +4.16%  "
$ws.Range("D22").Value = "'406.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("D23").Value = "'13.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +27.55%  "
$ws.Range("D24").Value = "'4.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'86.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("D27").Value = "'12.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.35%  "
$ws.Range("D28").Value = "'3.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.37%  "
$ws.Range("D29").Value = "'6.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").Value = "'8.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +23.03%  "
$ws.Range("D31").Value = "'9.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("D32").Value = "'32.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("D33").Value = "'697.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.50%  "
$ws.Range("D34").Value = "'12.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.34%  "
$ws.Range("D36").Value = "'64.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "'43.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("E38").Value = "  +15.88%  "
$ws.Range("D39").Value = "0.0₃0809"
$ws.Range("E39").Value = "  +9.10%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  +9.68%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'3.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.43%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.248.45"
$ws.Range("E43").Value = "  +11.20%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +19.49%  "
$ws.Range("D45").Value = "'3.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +37.88%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("D48").Value = "'9.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.70%  "
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("D50").Value = "'2.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.64%  "
$ws.Range("D51").Value = "'3.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.14%  "
